# Case 3.7 (380 kV) results update for res_line/pl_mw.xlsx, Sheet1
# Updates columns B,C,D,F,G,H,I,J,K,M,N for rows 2-25 (E, L remain 0, column A is the index, row 1 is header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 0.3752571294055258
    "C2" = 0.07041978470478227
    "D2" = 0.204324031431085
    "F2" = 1.801239737127617
    "G2" = 1.09076422471523
    "H2" = 1.087655301892994
    "I2" = 0.8761141638628835
    "J2" = 0.2871523237644809
    "K2" = 0.4267954579349578
    "M2" = 0.3168094642479886
    "N2" = 2.012933172946326
    "B3" = 0.3429065644524201
    "C3" = 0.06413048884360251
    "D3" = 0.2005809308349455
    "F3" = 1.800969117776837
    "G3" = 1.090735609142499
    "H3" = 1.092142826404995
    "I3" = 0.8804118492662631
    "J3" = 0.2853477372711168
    "K3" = 0.3896929908914899
    "M3" = 0.3039700844068989
    "N3" = 2.032261794613589
    "B4" = 0.3231634046853742
    "C4" = 0.06029348825607883
    "D4" = 0.198366495907436
    "F4" = 1.801679988168203
    "G4" = 1.0913229993713
    "H4" = 1.095336043500936
    "I4" = 0.8834550543796134
    "J4" = 0.2843783480104065
    "K4" = 0.3670515718673926
    "M4" = 0.2962399222303418
    "N4" = 2.044750757775269
    "B5" = 0.3151484902175241
    "C5" = 0.05873609904676869
    "D5" = 0.1974852640434364
    "F5" = 1.802190308250701
    "G5" = 1.091714433369788
    "H5" = 1.096747485925903
    "I5" = 0.884796884227125
    "J5" = 0.2840181981835102
    "K5" = 0.3578604670719017
    "M5" = 0.2931284983986515
    "N5" = 2.049996271594129
    "B6" = 0.31381947737529
    "C6" = 0.05847787202576171
    "D6" = 0.1973402171868059
    "F6" = 1.8022883752204
    "G6" = 1.091788613637192
    "H6" = 1.09698851179067
    "I6" = 0.8850258366520549
    "J6" = 0.2839605031957007
    "K6" = 0.3563364427633644
    "M6" = 0.2926141891325713
    "N6" = 2.050876716196541
    "B7" = 0.3230551884452382
    "C7" = 0.06027245954953742
    "D7" = 0.1983545254833103
    "F7" = 1.801685977054454
    "G7" = 1.09132766274476
    "H7" = 1.095354632489531
    "I7" = 0.8834727389824764
    "J7" = 0.2843733496275576
    "K7" = 0.3669274733946679
    "M7" = 0.2961978036813733
    "N7" = 2.044820868357469
    "B8" = 0.364077905023322
    "C8" = 0.06824612918039463
    "D8" = 0.2030160579753328
    "F8" = 1.800964453386854
    "G8" = 1.090628762565771
    "H8" = 1.089111767917046
    "I8" = 0.8775120740180107
    "J8" = 0.2865013416187594
    "K8" = 0.4139737370904015
    "M8" = 0.3123507139703605
    "N8" = 2.019468757601562
    "B9" = 0.4454659047803773
    "C9" = 0.08407801627046751
    "D9" = 0.2128194863532968
    "F9" = 1.806506322792956
    "G9" = 1.094062248966509
    "H9" = 1.080340889496028
    "I9" = 0.8690322095882692
    "J9" = 0.2917738607578286
    "K9" = 0.5073291908870488
    "M9" = 0.3452384622784734
    "N9" = 1.974681905070085
    "B10" = 0.5058273640812274
    "C10" = 0.09583008024745254
    "D10" = 0.2204225591644473
    "F10" = 1.814819471217035
    "G10" = 1.099521087373247
    "H10" = 1.076010131045706
    "I10" = 0.8647593071267465
    "J10" = 0.296318096673204
    "K10" = 0.5765803038505055
    "M10" = 0.3701371708984951
    "N10" = 1.944780964790436
    "B11" = 0.533408678087568
    "C11" = 0.1012028628595658
    "D11" = 0.223967742005442
    "F11" = 1.819522701196846
    "G11" = 1.102643819193673
    "H11" = 1.074498256156431
    "I11" = 0.8632407287303607
    "J11" = 0.2985310857249175
    "K11" = 0.6082274374857377
    "M11" = 0.3816236646031896
    "N11" = 1.931830412692545
    "B12" = 0.5438703598848065
    "C12" = 0.1032412348985474
    "D12" = 0.2253225777740369
    "F12" = 1.82143617812784
    "G12" = 1.103918379570061
    "H12" = 1.073991583037795
    "I12" = 0.8627268387184444
    "J12" = 0.2993900434568957
    "K12" = 0.6202319110706753
    "M12" = 0.3859962019708618
    "N12" = 1.927020101042682
    "B13" = 0.5416164891368567
    "C13" = 0.1028020656951298
    "D13" = 0.2250302416886854
    "F13" = 1.821018185972449
    "G13" = 1.103639784769143
    "H13" = 1.074097776625791
    "I13" = 0.8628347935549954
    "J13" = 0.2992041202136306
    "K13" = 0.6176456320831107
    "M13" = 0.385053483855593
    "N13" = 1.928051915285213
    "B14" = 0.5342690237326337
    "C14" = 0.1013704846205883
    "D14" = 0.2240789581095726
    "F14" = 1.81967746956839
    "G14" = 1.102746832696695
    "H14" = 1.074455252577394
    "I14" = 0.8631972247747086
    "J14" = 0.2986013329925044
    "K14" = 0.6092146462744665
    "M14" = 0.381982938757794
    "N14" = 1.931432785406145
    "B15" = 0.5297707221622829
    "C15" = 0.1004940960761189
    "D15" = 0.2234978754412964
    "F15" = 1.818873491540486
    "G15" = 1.102211864117706
    "H15" = 1.074682789939104
    "I15" = 0.8634271901895545
    "J15" = 0.2982348357705149
    "K15" = 0.6040530681967482
    "M15" = 0.3801051127510817
    "N15" = 1.933515881839757
    "B16" = 0.5040272936537633
    "C16" = 0.09547949160440794
    "D16" = 0.2201926063774238
    "F16" = 1.814530643462177
    "G16" = 1.099329885427991
    "H16" = 1.076118151768767
    "I16" = 0.8648671074500172
    "J16" = 0.2961764048214519
    "K16" = 0.5745149634350071
    "M16" = 0.3693897073512744
    "N16" = 1.945640432127703
    "B17" = 0.488265652793956
    "C17" = 0.0924100174398319
    "D17" = 0.2181870249971354
    "F17" = 1.812102438748624
    "G17" = 1.097725737076317
    "H17" = 1.077116023057201
    "I17" = 0.865859365545198
    "J17" = 0.2949509523934353
    "K17" = 0.5564310199091835
    "M17" = 0.3628570136190135
    "N17" = 1.953245388949263
    "B18" = 0.4792115353388056
    "C18" = 0.09064705262807138
    "D18" = 0.2170416153066412
    "F18" = 1.810792541997159
    "G18" = 1.096863258550442
    "H18" = 1.077733103867757
    "I18" = 0.8664701053944199
    "J18" = 0.2942598284102758
    "K18" = 0.5460432335104031
    "M18" = 0.3591146447622648
    "N18" = 1.957680878916037
    "B19" = 0.4761479643345581
    "C19" = 0.0900505758406922
    "D19" = 0.2166552013069776
    "F19" = 1.810363933092177
    "G19" = 1.096581572165931
    "H19" = 1.077949446232196
    "I19" = 0.8666837641149314
    "J19" = 0.2940281833390799
    "K19" = 0.5425284583256484
    "M19" = 0.3578501347159175
    "N19" = 1.95919318863892
    "B20" = 0.4899423126293243
    "C20" = 0.09273650803569922
    "D20" = 0.2183996799404326
    "F20" = 1.812351948137049
    "G20" = 1.097890271789282
    "H20" = 1.077005334345586
    "I20" = 0.8657495959809225
    "J20" = 0.2950799836621769
    "K20" = 0.558354679754018
    "M20" = 0.3635508720490748
    "N20" = 1.952429481271686
    "B21" = 0.5364266877392367
    "C21" = 0.1017908713220663
    "D21" = 0.2243580386260504
    "F21" = 1.820067675787826
    "G21" = 1.103006615421208
    "H21" = 1.074348466701224
    "I21" = 0.8630891098298079
    "J21" = 0.2987778178650586
    "K21" = 0.6116904808154402
    "M21" = 0.3828842132863102
    "N21" = 1.930437195955466
    "B22" = 0.5669072018078793
    "C22" = 0.10773065135038
    "D22" = 0.2283241265621569
    "F22" = 1.825882439147904
    "G22" = 1.106887020751145
    "H22" = 1.07299580033839
    "I22" = 0.8617068414468747
    "J22" = 0.3013166519239263
    "K22" = 0.6466671866144509
    "M22" = 0.39565279630974
    "N22" = 1.916610640268917
    "B23" = 0.5506301305551062
    "C23" = 0.1045584531497923
    "D23" = 0.2262007957374266
    "F23" = 1.822708360683592
    "G23" = 1.104766847502958
    "H23" = 1.073682646560641
    "I23" = 0.8624119559334318
    "J23" = 0.2999504635331789
    "K23" = 0.6279887276926104
    "M23" = 0.3888258312634605
    "N23" = 1.923940081583144
    "B24" = 0.4891842714527002
    "C24" = 0.09258889628796396
    "D24" = 0.2183035148816259
    "F24" = 1.812238876706616
    "G24" = 1.09781569947728
    "H24" = 1.077055241542709
    "I24" = 0.8657990973381544
    "J24" = 0.2950216068698097
    "K24" = 0.5574849653455942
    "M24" = 0.3632371370158651
    "N24" = 1.952798155959641
    "B25" = 0.4233483077464371
    "C25" = 0.07977401629659653
    "D25" = 0.2100968383195152
    "F25" = 1.804262281242742
    "G25" = 1.09261831247359
    "H25" = 1.082342328188346
    "I25" = 0.8709825957764181
    "J25" = 0.2902297704459045
    "K25" = 0.4819571968283753
    "M25" = 0.3362119925337055
    "N25" = 1.986270375419938
}

foreach ($ref in $newValues.Keys) {
    $ws.Range($ref).Value = $newValues[$ref]
}

Write-Output "Updated $($newValues.Count) cells on $($ws.Name)"
